# Apply the edits described by the commit "updated 4.0 files and mdl"
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "BAU Emissions" sheet: rename the model-run tag from "NoSettings" to
#    "test" in every label in column A (rows 4-280).
# ---------------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$rng = $wsBau.Range("A4:A280")
$rng.Replace(" : NoSettings", " : test", 2, 1, $false, $false, $true, $true)

# ---------------------------------------------------------------------------
# 2) "About" sheet: bump the "last updated" date (C1) from 3/18/2024 (45369)
#    to 4/5/2024 (45387).
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# ---------------------------------------------------------------------------
# 3) "BAU Emissions" sheet: refreshed model output for the
#    "Industrial Sector Energy Related Emissions before CCS[natural gas if,
#    iron and steel 241,CO2]" row (row 94), years 2032-2050 (columns M:AE).
# ---------------------------------------------------------------------------
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94").Value = 5005380
$wsBau.Range("R94").Value = 5005380
$wsBau.Range("S94").Value = 5005380
$wsBau.Range("T94").Value = 5005380
$wsBau.Range("U94").Value = 5005380
$wsBau.Range("V94").Value = 5005380
$wsBau.Range("W94").Value = 5005380
$wsBau.Range("X94").Value = 5005380
$wsBau.Range("Y94").Value = 5005380
$wsBau.Range("Z94").Value = 5005380
$wsBau.Range("AA94").Value = 5005380
$wsBau.Range("AB94").Value = 5005380
$wsBau.Range("AC94").Value = 5005380
$wsBau.Range("AD94").Value = 5005380
$wsBau.Range("AE94").Value = 5005380

# ---------------------------------------------------------------------------
# 4) Update the saved view/selection state:
#    - "BAU Emissions" scrolls/selects A30:AE280 instead of B283.
#    - "Current and Planned Capacity" is no longer the selected tab.
#    - "About" becomes the selected tab.
# ---------------------------------------------------------------------------
$wsBau.Activate()
$wsBau.Range("A30:AE280").Select()

$wsAbout.Activate()
